$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) to stay text even for values that look numeric,
# matching the workbook convention where D/E are plain text cells.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '67.314.44'
$ws.Range('E2').Value = '  +0.20%  '
$ws.Range('D3').Value = '3.107.82'
$ws.Range('E3').Value = '  -0.91%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').Value = '580.03'
$ws.Range('E5').Value = '  -0.01%  '
$ws.Range('D6').Value = '172.83'
$ws.Range('E6').Value = '  -1.26%  '
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('D8').Value = '3.101.76'
$ws.Range('E8').Value = '  -1.08%  '
$ws.Range('E9').Value = '  -0.98%  '
$ws.Range('E10').Value = '  -0.31%  '
$ws.Range('E11').Value = '  -1.70%  '
$ws.Range('D12').Value = '0.477'
$ws.Range('E12').Value = '  -1.47%  '
$ws.Range('D13').Value = '0.0000247'
$ws.Range('E13').Value = '  -1.45%  '
$ws.Range('D14').Value = '36.65'
$ws.Range('E14').Value = '  -2.10%  '
$ws.Range('E15').Value = '  -1.43%  '
$ws.Range('D16').Value = '3.626.62'
$ws.Range('E16').Value = '  -0.69%  '
$ws.Range('D17').Value = '67.322.77'
$ws.Range('E17').Value = '  +0.28%  '
$ws.Range('E18').Value = '  -1.08%  '
$ws.Range('D19').Value = '3.112.00'
$ws.Range('E19').Value = '  -0.84%  '
$ws.Range('D20').Value = '16.66'
$ws.Range('E20').Value = '  +2.93%  '
$ws.Range('D21').Value = '491.73'
$ws.Range('E21').Value = '  +0.75%  '
$ws.Range('D22').Value = '0.698'
$ws.Range('E22').Value = '  -2.82%  '
$ws.Range('D23').Value = '7.79'
$ws.Range('E23').Value = '  +1.88%  '
$ws.Range('D24').Value = '83.92'
$ws.Range('E24').Value = '  -0.56%  '
$ws.Range('D25').Value = '13.06'
$ws.Range('E25').Value = '  -1.14%  '
$ws.Range('D26').Value = '2.28'
$ws.Range('E26').Value = '  -2.47%  '
$ws.Range('D27').Value = '10.54'
$ws.Range('E27').Value = '  +4.27%  '
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('D29').Value = '7.85'
$ws.Range('E29').Value = '  -1.81%  '
$ws.Range('E30').Value = '  -2.74%  '
$ws.Range('E31').Value = '  -0.76%  '
$ws.Range('D32').Value = '28.31'
$ws.Range('E32').Value = '  -2.10%  '
$ws.Range('E33').Value = '  -1.40%  '
$ws.Range('D34').Value = '0.0₃0943'
$ws.Range('E34').Value = '  -5.60%  '
$ws.Range('E35').Value = '  +0.23%  '
$ws.Range('D36').Value = '5.79'
$ws.Range('E36').Value = '  -2.81%  '
$ws.Range('D37').Value = '0.972'
$ws.Range('E37').Value = '  -1.88%  '
$ws.Range('D38').Value = '46.64'
$ws.Range('E38').Value = '  -1.94%  '
$ws.Range('D39').Value = '2.03'
$ws.Range('E39').Value = '  -3.89%  '
$ws.Range('E40').Value = '  +0.58%  '
$ws.Range('D41').Value = '0.307'
$ws.Range('E41').Value = '  -2.09%  '
$ws.Range('D42').Value = '8.46'
$ws.Range('E42').Value = '  -2.60%  '
$ws.Range('D43').Value = '384.01'
$ws.Range('E43').Value = '  -0.35%  '
$ws.Range('D44').Value = '2.798.69'
$ws.Range('E44').Value = '  -1.94%  '
$ws.Range('D45').Value = '2.57'
$ws.Range('E45').Value = '  -8.77%  '
$ws.Range('E46').Value = '  -1.99%  '
$ws.Range('D47').Value = '135.32'
$ws.Range('E47').Value = '  -1.11%  '
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('D49').Value = '24.96'
$ws.Range('E49').Value = '  -0.27%  '
$ws.Range('D50').Value = '2.19'
$ws.Range('E50').Value = '  -1.48%  '
$ws.Range('E51').Value = '  -2.00%  '

# Drop the temporary text-number-format override so the cell styling
# matches the original (unstyled) cells -- only the values changed.
$ws.Range('D2:D51').ClearFormats()
